# Update "想去人数" (interest count, column F) values to the freshly-scraped
# numbers as published at commit 456a3b4 of the gh-pages data export.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 1059
$ws.Range("F6").Value = 2728
$ws.Range("F8").Value = 1337
$ws.Range("F9").Value = 945
$ws.Range("F10").Value = 643
$ws.Range("F11").Value = 956
$ws.Range("F12").Value = 1216
$ws.Range("F13").Value = 302
$ws.Range("F16").Value = 803
$ws.Range("F17").Value = 231
$ws.Range("F18").Value = 555
$ws.Range("F19").Value = 1148
$ws.Range("F21").Value = 673
$ws.Range("F23").Value = 237
$ws.Range("F24").Value = 327
$ws.Range("F25").Value = 320
$ws.Range("F26").Value = 702
$ws.Range("F27").Value = 656
$ws.Range("F28").Value = 6958
$ws.Range("F33").Value = 196
$ws.Range("F36").Value = 141
$ws.Range("F37").Value = 457
$ws.Range("F41").Value = 160
$ws.Range("F42").Value = 26
$ws.Range("F45").Value = 152
$ws.Range("F46").Value = 144

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 61
$ws.Range("F10").Value = 182
$ws.Range("F12").Value = 207
$ws.Range("F14").Value = 48
$ws.Range("F18").Value = 223
$ws.Range("F23").Value = 4

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 763

# Sheet 4: 全部类型 (All categories combined)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 1059
$ws.Range("F6").Value = 2728
$ws.Range("F7").Value = 1337
$ws.Range("F8").Value = 945
$ws.Range("F9").Value = 643
$ws.Range("F10").Value = 956
$ws.Range("F11").Value = 1216
$ws.Range("F12").Value = 302
$ws.Range("F16").Value = 803
$ws.Range("F17").Value = 231
$ws.Range("F18").Value = 555
$ws.Range("F19").Value = 1148
$ws.Range("F22").Value = 61
$ws.Range("F23").Value = 673
$ws.Range("F25").Value = 237
$ws.Range("F26").Value = 320
$ws.Range("F28").Value = 656
$ws.Range("F29").Value = 6958
$ws.Range("F30").Value = 207
$ws.Range("F32").Value = 196
$ws.Range("F34").Value = 48
$ws.Range("F35").Value = 48
$ws.Range("F39").Value = 160
$ws.Range("F40").Value = 26
$ws.Range("F44").Value = 152
$ws.Range("F48").Value = 4
